$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, column index, new text value.
# Column 4 (D/"Price") holds values that often look numeric (e.g. "0.999"),
# so those are written through a temporary Text ("@") number format to
# keep them as strings (preserving things like trailing zeros), and the
# cells original style is restored afterwards so no formatting changes
# leak into the workbook.
$updates = @(
    ,@(2, 4, '59.153.91')
    ,@(2, 5, '  -5.94%  ')
    ,@(3, 4, '2.440.53')
    ,@(3, 5, '  -8.68%  ')
    ,@(4, 4, '0.999')
    ,@(4, 5, '  -0.13%  ')
    ,@(5, 4, '531.22')
    ,@(5, 5, '  -3.63%  ')
    ,@(6, 4, '146.59')
    ,@(6, 5, '  -7.33%  ')
    ,@(7, 4, '0.997')
    ,@(7, 5, '  -0.29%  ')
    ,@(8, 4, '0.568')
    ,@(8, 5, '  -4.06%  ')
    ,@(9, 4, '2.449.81')
    ,@(9, 5, '  -8.30%  ')
    ,@(10, 4, '0.0988')
    ,@(10, 5, '  -6.41%  ')
    ,@(11, 5, '  -2.52%  ')
    ,@(12, 4, '5.22')
    ,@(12, 5, '  -0.24%  ')
    ,@(13, 4, '0.349')
    ,@(13, 5, '  -4.88%  ')
    ,@(14, 4, '2.868.14')
    ,@(14, 5, '  -8.92%  ')
    ,@(15, 4, '23.72')
    ,@(15, 5, '  -9.23%  ')
    ,@(16, 4, '59.096.09')
    ,@(16, 5, '  -5.87%  ')
    ,@(17, 4, '0.0000137')
    ,@(17, 5, '  -6.59%  ')
    ,@(18, 4, '2.492.14')
    ,@(18, 5, '  -6.90%  ')
    ,@(19, 4, '11.05')
    ,@(19, 5, '  -6.57%  ')
    ,@(20, 4, '4.33')
    ,@(20, 5, '  -5.70%  ')
    ,@(21, 4, '321.69')
    ,@(21, 5, '  -6.51%  ')
    ,@(22, 5, '  -3.33%  ')
    ,@(23, 4, '5.68')
    ,@(23, 5, '  -9.48%  ')
    ,@(24, 4, '0.460')
    ,@(24, 5, '  -8.92%  ')
    ,@(25, 4, '60.07')
    ,@(25, 5, '  -4.64%  ')
    ,@(26, 5, '  -4.55%  ')
    ,@(27, 5, '  -2.40%  ')
    ,@(28, 4, '7.65')
    ,@(28, 5, '  -6.07%  ')
    ,@(29, 4, '6.70')
    ,@(29, 5, '  -6.74%  ')
    ,@(30, 5, '  -7.03%  ')
    ,@(31, 4, '1.25')
    ,@(31, 5, '  -8.10%  ')
    ,@(32, 4, '0.0₃0764')
    ,@(32, 5, '  -10.73%  ')
    ,@(33, 5, '  -0.08%  ')
    ,@(34, 4, '154.42')
    ,@(34, 5, '  -7.96%  ')
    ,@(35, 4, '1.37')
    ,@(35, 5, '  -7.30%  ')
    ,@(36, 4, '4.51')
    ,@(36, 5, '  -6.91%  ')
    ,@(37, 4, '18.28')
    ,@(37, 5, '  -6.17%  ')
    ,@(38, 5, '  -3.55%  ')
    ,@(39, 5, '  -8.12%  ')
    ,@(40, 4, '310.41')
    ,@(40, 5, '  -10.88%  ')
    ,@(41, 5, '  -4.23%  ')
    ,@(42, 2, 'SuiNetwork')
    ,@(42, 3, 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui')
    ,@(42, 4, '0.825')
    ,@(42, 5, '  -13.83%  ')
    ,@(43, 2, 'Filecoin')
    ,@(43, 3, 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil')
    ,@(43, 4, '3.68')
    ,@(43, 5, '  -7.64%  ')
    ,@(44, 4, '0.995')
    ,@(44, 5, '  -0.35%  ')
    ,@(45, 4, '10.72')
    ,@(45, 5, '  -2.85%  ')
    ,@(46, 4, '0.0935')
    ,@(46, 5, '  -3.83%  ')
    ,@(47, 4, '0.578')
    ,@(47, 5, '  -6.15%  ')
    ,@(48, 4, '0.0522')
    ,@(48, 5, '  -7.22%  ')
    ,@(49, 4, '0.0227')
    ,@(49, 5, '  -5.68%  ')
    ,@(50, 4, '18.71')
    ,@(50, 5, '  -9.76%  ')
    ,@(51, 2, 'Aave')
    ,@(51, 3, 'https://coinranking.com/coin/ixgUfzmLR+aave-aave')
    ,@(51, 4, '120.41')
    ,@(51, 5, '  -6.42%  ')
)

foreach ($u in $updates) {
    $row = $u[0]
    $col = $u[1]
    $val = $u[2]
    $cell = $ws.Cells.Item($row, $col)
    if ($col -eq 4) {
        $origStyle = $cell.Style
        $cell.NumberFormat = "@"
        $cell.Value = $val
        $cell.Style = $origStyle
    } else {
        $cell.Value = $val
    }
}
